# "New Script - IPA" -- append two new test-case rows (IPA0003 / IPA0004)
# to the "Test Cases" sheet, matching the formatting of the row above them,
# plus a JIRA hyperlink for the new OPQA-4221 reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 9: IPA0003 ---------------------------------------------------
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("A9:E9").PasteSpecial(-4122) | Out-Null

# --- Row 10: IPA0004 ---------------------------------------------------
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("A10:E10").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Values -- write in the same order the new unique shared strings were
# originally authored in, so the shared-string table grows the same way.
$ws.Range("B9").Value = "OPQA-4241||OPQA-4245"
$ws.Range("C9").Value = 'Verify that user can skip the linking by clicking on "Not now button" on the modal "Already have an account? .. || Verify that once the user skips linking then user will not be prompted to link again.'
$ws.Range("A9").Value = "IPA0003"
$ws.Range("D9").Value = "Y"

$ws.Range("A10").Value = "IPA0004"
$ws.Range("C10").Value = 'Verify that error message " Incorrect password. Please try again."should be displayed when user enters incorrect password for existing steam account.|| Verify that when user''s account is locked due to 10 invalid authentications of existing password,user becomes locked, the user is signed out'
$ws.Range("B10").Value = "OPQA-4221 || OPQA-4225"
$ws.Range("D10").Value = "Y"

# Row heights (wrapped description text spans multiple lines).
$ws.Rows(9).RowHeight = 45
$ws.Rows(10).RowHeight = 75

# Hyperlink on the new OPQA-4221 JIRA reference cell.
$ws.Hyperlinks.Add($ws.Range("B10"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4221") | Out-Null

# Hyperlinks.Add auto-applies the built-in blue/underline "Hyperlink" style;
# restore the plain bordered look the rest of the table uses (the source
# workbook keeps hyperlinked cells in the ordinary table format).
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("A10:E10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the new last row selected, as in the authored workbook.
$ws.Range("A10:E10").Select() | Out-Null
